$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original hyperlink data: "<old row>|<url>" (row numbers as in the sheet
# BEFORE the new header row is inserted).
$hyperlinkData = @(
    "1|http://www.maavianrecords.com/review-list/black-bellied-whistling-duck",
    "2|http://www.maavianrecords.com/review-list/pink-footed-goose",
    "3|http://www.maavianrecords.com/review-list/barnacle-goose",
    "4|http://www.maavianrecords.com/review-list/tundra-swan",
    "5|http://www.maavianrecords.com/review-list/tufted-duck",
    "6|http://www.maavianrecords.com/review-list/eared-grebe",
    "7|http://www.maavianrecords.com/review-list/white-winged-dove",
    "8|http://www.maavianrecords.com/review-list/rufous-hummingbird",
    "9|http://www.maavianrecords.com/review-list/yellow-rail",
    "10|http://www.maavianrecords.com/review-list/purple-gallinule",
    "11|http://www.maavianrecords.com/review-list/northern-lapwing",
    "12|http://www.maavianrecords.com/review-list/wilsons-plover",
    "13|http://www.maavianrecords.com/review-list/bar-tailed-godwit",
    "14|http://www.maavianrecords.com/review-list/ruff",
    "15|http://www.maavianrecords.com/review-list/curlew-sandpiper",
    "16|http://www.maavianrecords.com/review-list/great-skua",
    "17|http://www.maavianrecords.com/review-list/south-polar-skua",
    "18|http://www.maavianrecords.com/review-list/ivory-gull",
    "19|http://www.maavianrecords.com/review-list/franklins-gull",
    "20|http://www.maavianrecords.com/review-list/mew-gull",
    "21|http://www.maavianrecords.com/review-list/sooty-tern",
    "22|http://www.maavianrecords.com/review-list/bridled-tern",
    "23|http://www.maavianrecords.com/review-list/gull-billed-tern",
    "24|http://www.maavianrecords.com/review-list/pacific-loon",
    "25|http://www.maavianrecords.com/review-list/brown-booby",
    "26|http://www.maavianrecords.com/review-list/american-white-pelican",
    "27|http://www.maavianrecords.com/review-list/brown-pelican",
    "28|http://www.maavianrecords.com/review-list/white-ibis",
    "29|http://www.maavianrecords.com/review-list/white-faced-ibis",
    "30|http://www.maavianrecords.com/review-list/swallow-tailed-kite",
    "31|http://www.maavianrecords.com/review-list/white-tailed-kite",
    "32|http://www.maavianrecords.com/review-list/mississippi-kite",
    "33|http://www.maavianrecords.com/review-list/black-backed-woodpecker",
    "34|http://www.maavianrecords.com/review-list/says-phoebe",
    "35|http://www.maavianrecords.com/review-list/ash-throated-flycatcher",
    "36|http://www.maavianrecords.com/review-list/scissor-tailed-flycatcher",
    "37|http://www.maavianrecords.com/review-list/fork-tailed-flycatcher",
    "38|http://www.maavianrecords.com/review-list/loggerhead-shrike",
    "39|http://www.maavianrecords.com/review-list/bells-vireo",
    "40|http://www.maavianrecords.com/review-list/cave-swallow",
    "41|http://www.maavianrecords.com/review-list/boreal-chickadee",
    "42|http://www.maavianrecords.com/review-list/northern-wheatear",
    "43|http://www.maavianrecords.com/review-list/mountain-bluebird",
    "44|http://www.maavianrecords.com/review-list/townsends-solitaire",
    "45|http://www.maavianrecords.com/review-list/varied-thrush",
    "46|http://www.maavianrecords.com/review-list/black-throated-gray-warbler",
    "47|http://www.maavianrecords.com/review-list/townsends-warbler",
    "48|http://www.maavianrecords.com/review-list/green-tailed-towhee",
    "50|http://www.maavianrecords.com/review-list/le-contes-sparrow",
    "51|http://www.maavianrecords.com/review-list/harriss-sparrow",
    "52|http://www.maavianrecords.com/review-list/western-tanager",
    "53|http://www.maavianrecords.com/review-list/black-headed-grosbeak",
    "54|http://www.maavianrecords.com/review-list/painted-bunting"
)

# Insert a new row at the top; existing rows (and their values/styles) shift
# down by one automatically.
$ws.Rows.Item(1).Insert()

# New header cell for the shifted-down species list.
$ws.Range("A1").Value = "Species"

# Capture the (already shifted) cell text for every row that carries a
# hyperlink, so it can be restored after the hyperlink is recreated.
$savedValues = @{}
foreach ($entry in $hyperlinkData) {
    $oldRow = [int]($entry.Split("|", 2)[0])
    $newRow = $oldRow + 1
    $savedValues[$newRow] = $ws.Cells.Item($newRow, 1).Value2
}

# The row insert does not relocate the worksheet's hyperlinks in this
# runtime, so clear all existing hyperlinks and re-create them one row
# lower, pointing at the same urls as before. Adding a hyperlink with
# TextToDisplay overwrites the cell's shown text as a side effect, so the
# original species name and Hyperlink style are written back immediately
# afterward.
$ws.Cells.Hyperlinks.Delete()

foreach ($entry in $hyperlinkData) {
    $parts = $entry.Split("|", 2)
    $oldRow = [int]$parts[0]
    $url = $parts[1]
    $newRow = $oldRow + 1
    $cell = $ws.Cells.Item($newRow, 1)
    $ws.Hyperlinks.Add($cell, $url, "", "", $url) | Out-Null
    $cell.Value = $savedValues[$newRow]
    $cell.Style = "Hyperlink"
}
